$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "image.path" column before the existing "description" column (C),
# pushing description/tool.type/external.data.support one column to the right.
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Header + image path values for the new column C
$ws.Cells.Item(1,3).Value = "image.path"
$ws.Cells.Item(2,3).Value = "shiny-logo.png"
$ws.Cells.Item(3,3).Value = "plotly-logo.png"
$ws.Cells.Item(4,3).Value = "r-logo.png"
$ws.Cells.Item(5,3).Value = "dash-logo.svg"
$ws.Cells.Item(6,3).Value = "vega-lite-logo.png"

# Updated descriptions (column D, after the column insert) for Vega-Lite and Dash
$ws.Cells.Item(6,4).Value = 'Vega-Lite provides a high-level grammar of interactive graphics, allowing users to specify "charts as data" in well designed JSON format. OxShef are currently investigating the reproducability of a Vega-Lite driven dataviz workflow, which may lead to a dedicated site about this tool in the future.'
$ws.Cells.Item(5,4).Value = 'Dash allows Python users to build rich interactive web applications and visualisations through a combination of different technologies, including React and Flask. OxShef are currently investigating the reproducability of a Dash-driven dataviz visualisation workflow, which may lead to a dedicated site about this tool in the future.'

# Row 6 (Vega-Lite) height changed from 60 to 75
$ws.Rows.Item(6).RowHeight = 75

# New row 7: Jupyter
$ws.Cells.Item(7,1).Value = "Jupyter"
$ws.Cells.Item(7,2).Value = "http://jupyter.org/"
$ws.Cells.Item(7,3).Value = "jupyter-logo.png"
$ws.Cells.Item(7,4).Value = 'Jupyter (the spiritual successor to iPython notebooks) is a powerful tool for creating rich documents incorporating code, data and visualisation outputs. Jupyter notebooks allow code written in Python, R and more to be combined together easily. Oxshef are currently developing a site dedicated to using this tool in a reproducible dataviz workflow.'
$ws.Cells.Item(7,4).WrapText = $true
$ws.Cells.Item(7,5).Value = "Scripting"
$ws.Cells.Item(7,6).Value = $true
$ws.Rows.Item(7).RowHeight = 90

$null = $ws.Range("F7").Select()

Write-Output "done"
